$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the template/filler rows that belonged to the first "Brake Master
# Cylinder" assembly block (old rows 4-8) and to the tail of the "Balance
# bar" assembly block (old rows 12-15), leaving a 6-row sheet.
$ws.Range("A4:G8").EntireRow.Delete()
$ws.Range("A7:G10").EntireRow.Delete()

# Row 3: first real part entry under "Brake Master Cylinder"
$ws.Cells.Item(3, 3).Value = "Master Cylinder"
$ws.Cells.Item(3, 4).Value = "b"
$ws.Cells.Item(3, 6).Value = 2

# Row 5: first real part entry under "Balance bar"
$ws.Cells.Item(5, 3).Value = "Balance bar"
$ws.Cells.Item(5, 4).Value = "b"
$ws.Cells.Item(5, 6).Value = 2

# Row 6: second real part entry under "Balance bar"
$ws.Cells.Item(6, 3).Value = "Master Cylinder support"
$ws.Cells.Item(6, 4).Value = "m"
$ws.Cells.Item(6, 5).Value = "Threaded aluminum part which support the bottom of the master cylinder"
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = "BR_04002"

# The long wrapped description makes row 6 taller (two lines instead of one).
$ws.Rows.Item(6).RowHeight = 27

$ws.Range("D10").Select()
